$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add I0 in I1 and IF in J1, matching the style/format used by existing header cells (e.g. H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows 2-7: add value 9 in columns I and J
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 9).Value = 9
    $ws.Cells.Item($r, 10).Value = 9
}
